$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 24-40 of column D currently have no cell/style at all; rows 21-23
# already carry the "s=1" style (Arial, left/center, indent 1) on empty
# cells. Copy that style down to D24:D40 first so every cell D21:D40
# ends up with the same formatting, then fill in the pass/fail values.
$ws.Range("D21").Copy()
$ws.Range("D24:D40").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("D21").Value = "pass"
$ws.Range("D22").Value = "paas"
$ws.Range("D23").Value = "pass"
$ws.Range("D24").Value = "pass"
$ws.Range("D25").Value = "pass"
$ws.Range("D26").Value = "fail"
$ws.Range("D27").Value = "fail"
$ws.Range("D28").Value = "pass"
$ws.Range("D29").Value = "pass"
$ws.Range("D30").Value = "pass"
$ws.Range("D31").Value = "pass"
$ws.Range("D32").Value = "fail"
$ws.Range("D33").Value = "pass"
$ws.Range("D34").Value = "pass"
$ws.Range("D35").Value = "pass"
$ws.Range("D36").Value = "pass"
$ws.Range("D37").Value = "pass"
$ws.Range("D38").Value = "pass"
$ws.Range("D39").Value = "fail"
$ws.Range("D40").Value = "pass"

# Move the active selection to B32, matching the saved view state.
$ws.Range("B32").Select()
